# Reorganize test data: append 4 new claim rows (22-25) describing
# multi-vendor / services purchase orders for Northstar Technologies Inc,
# Valley Data Systems and Harbor Communications.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: Northstar Technologies Inc / Hardware Tracker PO ---
$ws.Range("A22").Value = "V0021"
$ws.Range("B22").Value = "Northstar Technologies Inc"
$ws.Range("E22").Value = "PO-Hardware-Tracker.xlsx"
$ws.Range("H22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H22").Value = "2024-05-11"
$ws.Range("I22").Value = 932805
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 95612.50999999999
$ws.Range("L22").Value = 932805
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = "IT Equipment (Multi-vendor PO)"
$ws.Range("O22").Value = "UPS Battery Backup, Server Rack, System Integration, Storage Array"
$ws.Range("Q22").Value = "Hardware Tracker, Multi-vendor"
$ws.Range("R22").Value = "Seattle, WA"
$ws.Range("S22").Value = 10.25
$ws.Range("T22").Value = 95612.50999999999
$ws.Range("U22").Value = 0

# --- Row 23: Northstar Technologies Inc / Software Summary PO ---
$ws.Range("A23").Value = "V0022"
$ws.Range("B23").Value = "Northstar Technologies Inc"
$ws.Range("E23").Value = "PO-Software-Summary.xlsx"
$ws.Range("H23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H23").Value = "2024-04-28"
$ws.Range("I23").Value = 442996
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 45407.09
$ws.Range("L23").Value = 442996
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = "Software License (Multi-vendor PO)"
$ws.Range("O23").Value = "Development Tools, Enterprise License, Firewall, Security Suite, Managed Services"
$ws.Range("Q23").Value = "Software Summary, Multi-vendor"
$ws.Range("R23").Value = "Seattle, WA"
$ws.Range("S23").Value = 10.25
$ws.Range("T23").Value = 45407.09
$ws.Range("U23").Value = 0

# --- Row 24: Valley Data Systems / Consulting Approval Email ---
$ws.Range("A24").Value = "V0023"
$ws.Range("B24").Value = "Valley Data Systems"
$ws.Range("E24").Value = "PO-Approval-Email.eml"
$ws.Range("G24").Value = "PO-4900348121"
$ws.Range("H24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H24").Value = "2023-02-14"
$ws.Range("I24").Value = 1013687
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 103902.92
$ws.Range("L24").Value = 1013687
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = "Professional Services"
$ws.Range("O24").Value = "Training Program (7), Technical Consulting (8), System Integration Services (9)"
$ws.Range("Q24").Value = "Training, Consulting, System Integration"
$ws.Range("R24").Value = "Seattle, WA"
$ws.Range("S24").Value = 10.25
$ws.Range("T24").Value = 103902.92
$ws.Range("U24").Value = 0

# --- Row 25: Harbor Communications / Consulting Quotation ---
$ws.Range("A25").Value = "V0024"
$ws.Range("B25").Value = "Harbor Communications"
$ws.Range("E25").Value = "PO-Consulting-Quotation.docx"
$ws.Range("G25").Value = "Q-73483"
$ws.Range("H25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H25").Value = "2024-07-05"
$ws.Range("I25").Value = 569214
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 58344.44
$ws.Range("L25").Value = 569214
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = "Consulting Services"
$ws.Range("O25").Value = "Telecom services implementation"
$ws.Range("Q25").Value = "Consulting Services Proposal"
$ws.Range("R25").Value = "Seattle, WA"
$ws.Range("S25").Value = 10.25
$ws.Range("T25").Value = 58344.44
$ws.Range("U25").Value = 0
